# Applies the "Fixed images and files" commit to language_definition.xlsx
#
# Summary of semantic changes (see commit diff):
#  - Sheet "Language Definition" (Custom Concept table):
#      Row 30 (KG-25-27 StartDate): "Date - Schema.org" -> "Date - W3.org"
#      Row 31 (KG-25-28 EndDate):   rich-text run "Schema.org" -> "W3.org"
#                                   (keeps the blue/underlined hyperlink-style run)
#      Row 33 (KG-25-29 MicroClimate): Equivalent Concept "Custom" -> "SweetOntology "
#      Row 39 (KG-25-34): Custom Concept "Region" -> "City"
#                          Gloss "Code of Region." -> "City of Trentino Region"
#      Row 40 (KG-25-35): Gloss "Name of Region." -> "Name of the city"
#  - Sheet "Knowledge - Data link":
#      Row 39 (KG-25-34): Custom Concept "Region" -> "City" (same concept rename)

$wb = $excel.ActiveWorkbook

$wsLang = $wb.Worksheets.Item("Language Definition")
$wsLink = $wb.Worksheets.Item("Knowledge - Data link")

# --- Language Definition sheet -------------------------------------------------

# Row 30 - StartDate: plain-text "Date - Schema.org" -> "Date - W3.org"
$wsLang.Cells.Item(30, 3).Value = "Date - W3.org"

# Row 31 - EndDate: rich text run "Schema.org" -> "W3.org", keep the
# existing hyperlink-style formatting (blue FF1155CC + underline) on the
# "W3.org" portion, plain formatting on the "Date - " portion.
$endDateCell = $wsLang.Cells.Item(31, 3)
$endDateCell.Value = "Date - W3.org"
$endDateRun = $endDateCell.Characters(8, 6)
$endDateRun.Font.Color = 13391121
$endDateRun.Font.Underline = $true

# Row 33 - MicroClimate: Equivalent Concept "Custom" -> "SweetOntology "
$wsLang.Cells.Item(33, 3).Value = "SweetOntology "

# Row 39 - Region -> City
$wsLang.Cells.Item(39, 2).Value = "City"
$wsLang.Cells.Item(39, 4).Value = "City of Trentino Region"

# Row 40 - Name of Region. -> Name of the city
$wsLang.Cells.Item(40, 4).Value = "Name of the city"

# --- Knowledge - Data link sheet ------------------------------------------------

# Row 39 - Region -> City (same Custom Concept rename reflected here too)
$wsLink.Cells.Item(39, 2).Value = "City"
